# #deleted unused variables in Molten Salt Tower Parasitics UI page
#
# Appends 7 new "Deleted variable" rows (36-42) to the "SAM Variable Changes"
# sheet, documenting unused variables removed from the Molten Salt Tower
# Parasitics UI page.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Duplicate the last existing data row (35) seven times, inserting the
# copies directly below it, so the new rows inherit the same per-cell
# styles (fills/alignment) as the rest of the "Deleted variable" block.
for ($i = 0; $i -lt 7; $i++) {
    $ws.Rows(35).Copy() | Out-Null
    $ws.Rows(36).Insert() | Out-Null
}

# Rows 38-42 use the plain (unshaded) look for column F, matching the
# style already used by F31 rather than the copied F35 style.
$ws.Range("F31").Copy() | Out-Null
$ws.Range("F38:F42").PasteSpecial(-4122) | Out-Null

# Fill in the variable-specific text for each new row. Columns A, B, E, G
# and H already carry the correct values from the row-35 copy
# (Deleted variable / number / Molten Salt Tower Parasitics / N/A / Ty).
$ws.Range("C36").Value = "P_storage_pump"
$ws.Range("C37").Value = "storage_bypass"

$ws.Range("C38").Value = "recirc_source"
$ws.Range("F38").Value = "not used"

$ws.Range("C39").Value = "recirc_htf_eff"
$ws.Range("F39").Value = "not used"

$ws.Range("C40").Value = "flow_from_storage"
$ws.Range("F40").Value = "not used"

$ws.Range("C41").Value = "P_hot_tank"
$ws.Range("F41").Value = "not used"

$ws.Range("C42").Value = "csp.pt.par.bop_c1"
$ws.Range("F42").Value = "not used"

# Match the saved view state: scrolled down a bit, selection on A43.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A43").Select() | Out-Null
